$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"21.83470933333334"
$ws.Cells.Item(2, 8).Value = [double]"65.50412800000001"
$ws.Cells.Item(2, 9).Value = [double]"0.1994658397831471"
$ws.Cells.Item(2, 10).Value = [double]"0.1994658397831471"
$ws.Cells.Item(2, 11).Value = [double]"3"
$ws.Cells.Item(2, 12).Value = [double]"1"
$ws.Cells.Item(2, 13).Value = [double]"3.825035"
$ws.Cells.Item(2, 14).Value = [double]"11.475105"
$ws.Cells.Item(2, 15).Value = [double]"0.03111562857396839"
$ws.Cells.Item(2, 16).Value = [double]"0.03111562857396839"
$ws.Cells.Item(2, 17).Value = [double]"83.51852741482668"
$ws.Cells.Item(2, 18).Value = [double]"751.66674673344"
$ws.Cells.Item(2, 19).Value = [double]"0.006206504983887093"
$ws.Cells.Item(2, 20).Value = [double]"0.006206504983887093"
$ws.Cells.Item(3, 7).Value = [double]"21.83470933333334"
$ws.Cells.Item(3, 8).Value = [double]"65.50412800000001"
$ws.Cells.Item(3, 9).Value = [double]"0.1994658397831471"
$ws.Cells.Item(3, 10).Value = [double]"0.1994658397831471"
$ws.Cells.Item(3, 15).Value = [double]"0.4709815605157605"
$ws.Cells.Item(3, 16).Value = [double]"0.4709815605157605"
$ws.Cells.Item(3, 17).Value = [double]"1264.177783852388"
$ws.Cells.Item(3, 18).Value = [double]"11377.60005467149"
$ws.Cells.Item(3, 19).Value = [double]"0.0939447324906533"
$ws.Cells.Item(3, 20).Value = [double]"0.09394473249065328"
$ws.Cells.Item(4, 7).Value = [double]"21.83470933333334"
$ws.Cells.Item(4, 8).Value = [double]"65.50412800000001"
$ws.Cells.Item(4, 9).Value = [double]"0.1994658397831471"
$ws.Cells.Item(4, 10).Value = [double]"0.1994658397831471"
$ws.Cells.Item(4, 13).Value = [double]"61.10114166666667"
$ws.Cells.Item(4, 14).Value = [double]"183.303425"
$ws.Cells.Item(4, 15).Value = [double]"0.4970413158429724"
$ws.Cells.Item(4, 16).Value = [double]"0.4970413158429724"
$ws.Cells.Item(4, 17).Value = [double]"1334.125668226489"
$ws.Cells.Item(4, 18).Value = [double]"12007.1310140384"
$ws.Cells.Item(4, 19).Value = [double]"0.09914276347153897"
$ws.Cells.Item(4, 20).Value = [double]"0.09914276347153896"
$ws.Cells.Item(5, 7).Value = [double]"21.83470933333334"
$ws.Cells.Item(5, 8).Value = [double]"65.50412800000001"
$ws.Cells.Item(5, 9).Value = [double]"0.1994658397831471"
$ws.Cells.Item(5, 10).Value = [double]"0.1994658397831471"
$ws.Cells.Item(5, 13).Value = [double]"0.1059033333333333"
$ws.Cells.Item(5, 14).Value = [double]"0.31771"
$ws.Cells.Item(5, 15).Value = [double]"0.0008614950672987739"
$ws.Cells.Item(5, 16).Value = [double]"0.0008614950672987739"
$ws.Cells.Item(5, 17).Value = [double]"2.312368500764445"
$ws.Cells.Item(5, 18).Value = [double]"20.81131650688"
$ws.Cells.Item(5, 19).Value = [double]"0.0001718388370677888"
$ws.Cells.Item(5, 20).Value = [double]"0.0001718388370677888"
$ws.Cells.Item(6, 8).Value = [double]"92.79671999999999"
$ws.Cells.Item(6, 9).Value = [double]"0.282574186529459"
$ws.Cells.Item(6, 10).Value = [double]"0.282574186529459"
$ws.Cells.Item(6, 11).Value = [double]"3"
$ws.Cells.Item(6, 12).Value = [double]"1"
$ws.Cells.Item(6, 13).Value = [double]"3.825035"
$ws.Cells.Item(6, 14).Value = [double]"11.475105"
$ws.Cells.Item(6, 15).Value = [double]"0.03111562857396839"
$ws.Cells.Item(6, 16).Value = [double]"0.03111562857396839"
$ws.Cells.Item(6, 17).Value = [double]"118.3169006284"
$ws.Cells.Item(6, 18).Value = [double]"1064.8521056556"
$ws.Cells.Item(6, 19).Value = [double]"0.008792473432641909"
$ws.Cells.Item(6, 20).Value = [double]"0.008792473432641907"
$ws.Cells.Item(7, 8).Value = [double]"92.79671999999999"
$ws.Cells.Item(7, 9).Value = [double]"0.282574186529459"
$ws.Cells.Item(7, 10).Value = [double]"0.282574186529459"
$ws.Cells.Item(7, 15).Value = [double]"0.4709815605157605"
$ws.Cells.Item(7, 16).Value = [double]"0.4709815605157605"
$ws.Cells.Item(7, 17).Value = [double]"1790.903190687013"
$ws.Cells.Item(7, 19).Value = [double]"0.1330872313331162"
$ws.Cells.Item(7, 20).Value = [double]"0.1330872313331162"
$ws.Cells.Item(8, 8).Value = [double]"92.79671999999999"
$ws.Cells.Item(8, 9).Value = [double]"0.282574186529459"
$ws.Cells.Item(8, 10).Value = [double]"0.282574186529459"
$ws.Cells.Item(8, 13).Value = [double]"61.10114166666667"
$ws.Cells.Item(8, 14).Value = [double]"183.303425"
$ws.Cells.Item(8, 15).Value = [double]"0.4970413158429724"
$ws.Cells.Item(8, 16).Value = [double]"0.4970413158429724"
$ws.Cells.Item(8, 17).Value = [double]"1889.995178307333"
$ws.Cells.Item(8, 18).Value = [double]"17009.956604766"
$ws.Cells.Item(8, 19).Value = [double]"0.1404510454958598"
$ws.Cells.Item(8, 20).Value = [double]"0.1404510454958598"
$ws.Cells.Item(9, 8).Value = [double]"92.79671999999999"
$ws.Cells.Item(9, 9).Value = [double]"0.282574186529459"
$ws.Cells.Item(9, 10).Value = [double]"0.282574186529459"
$ws.Cells.Item(9, 13).Value = [double]"0.1059033333333333"
$ws.Cells.Item(9, 14).Value = [double]"0.31771"
$ws.Cells.Item(9, 15).Value = [double]"0.0008614950672987739"
$ws.Cells.Item(9, 16).Value = [double]"0.0008614950672987739"
$ws.Cells.Item(9, 17).Value = [double]"3.275827323466666"
$ws.Cells.Item(9, 18).Value = [double]"29.4824459112"
$ws.Cells.Item(9, 19).Value = [double]"0.0002434362678410926"
$ws.Cells.Item(9, 20).Value = [double]"0.0002434362678410925"
$ws.Cells.Item(10, 7).Value = [double]"47.70664233333334"
$ws.Cells.Item(10, 8).Value = [double]"143.119927"
$ws.Cells.Item(10, 9).Value = [double]"0.435812784634851"
$ws.Cells.Item(10, 10).Value = [double]"0.435812784634851"
$ws.Cells.Item(10, 11).Value = [double]"3"
$ws.Cells.Item(10, 12).Value = [double]"1"
$ws.Cells.Item(10, 13).Value = [double]"3.825035"
$ws.Cells.Item(10, 14).Value = [double]"11.475105"
$ws.Cells.Item(10, 15).Value = [double]"0.03111562857396839"
$ws.Cells.Item(10, 16).Value = [double]"0.03111562857396839"
$ws.Cells.Item(10, 17).Value = [double]"182.4795766574817"
$ws.Cells.Item(10, 18).Value = [double]"1642.316189917335"
$ws.Cells.Item(10, 19).Value = [double]"0.0135605887344849"
$ws.Cells.Item(10, 20).Value = [double]"0.0135605887344849"
$ws.Cells.Item(11, 7).Value = [double]"47.70664233333334"
$ws.Cells.Item(11, 8).Value = [double]"143.119927"
$ws.Cells.Item(11, 9).Value = [double]"0.435812784634851"
$ws.Cells.Item(11, 10).Value = [double]"0.435812784634851"
$ws.Cells.Item(11, 15).Value = [double]"0.4709815605157605"
$ws.Cells.Item(11, 16).Value = [double]"0.4709815605157605"
$ws.Cells.Item(11, 17).Value = [double]"2762.101224215603"
$ws.Cells.Item(11, 18).Value = [double]"24858.91101794042"
$ws.Cells.Item(11, 19).Value = [double]"0.2052597854000412"
$ws.Cells.Item(11, 20).Value = [double]"0.2052597854000411"
$ws.Cells.Item(12, 7).Value = [double]"47.70664233333334"
$ws.Cells.Item(12, 8).Value = [double]"143.119927"
$ws.Cells.Item(12, 9).Value = [double]"0.435812784634851"
$ws.Cells.Item(12, 10).Value = [double]"0.435812784634851"
$ws.Cells.Item(12, 13).Value = [double]"61.10114166666667"
$ws.Cells.Item(12, 14).Value = [double]"183.303425"
$ws.Cells.Item(12, 15).Value = [double]"0.4970413158429724"
$ws.Cells.Item(12, 16).Value = [double]"0.4970413158429724"
$ws.Cells.Item(12, 17).Value = [double]"2914.930311649998"
$ws.Cells.Item(12, 18).Value = [double]"26234.37280484998"
$ws.Cells.Item(12, 19).Value = [double]"0.2166169599360963"
$ws.Cells.Item(12, 20).Value = [double]"0.2166169599360963"
$ws.Cells.Item(13, 7).Value = [double]"47.70664233333334"
$ws.Cells.Item(13, 8).Value = [double]"143.119927"
$ws.Cells.Item(13, 9).Value = [double]"0.435812784634851"
$ws.Cells.Item(13, 10).Value = [double]"0.435812784634851"
$ws.Cells.Item(13, 13).Value = [double]"0.1059033333333333"
$ws.Cells.Item(13, 14).Value = [double]"0.31771"
$ws.Cells.Item(13, 15).Value = [double]"0.0008614950672987739"
$ws.Cells.Item(13, 16).Value = [double]"0.0008614950672987739"
$ws.Cells.Item(13, 17).Value = [double]"5.052292445241112"
$ws.Cells.Item(13, 18).Value = [double]"45.47063200717"
$ws.Cells.Item(13, 19).Value = [double]"0.000375450564228667"
$ws.Cells.Item(13, 20).Value = [double]"0.000375450564228667"
$ws.Cells.Item(14, 7).Value = [double]"8.992316666666666"
$ws.Cells.Item(14, 8).Value = [double]"26.97695"
$ws.Cells.Item(14, 9).Value = [double]"0.08214718905254291"
$ws.Cells.Item(14, 10).Value = [double]"0.08214718905254291"
$ws.Cells.Item(14, 11).Value = [double]"3"
$ws.Cells.Item(14, 12).Value = [double]"1"
$ws.Cells.Item(14, 13).Value = [double]"3.825035"
$ws.Cells.Item(14, 14).Value = [double]"11.475105"
$ws.Cells.Item(14, 15).Value = [double]"0.03111562857396839"
$ws.Cells.Item(14, 16).Value = [double]"0.03111562857396839"
$ws.Cells.Item(14, 17).Value = [double]"34.39592598108333"
$ws.Cells.Item(14, 18).Value = [double]"309.56333382975"
$ws.Cells.Item(14, 19).Value = [double]"0.002556061422954487"
$ws.Cells.Item(14, 20).Value = [double]"0.002556061422954488"
$ws.Cells.Item(15, 7).Value = [double]"8.992316666666666"
$ws.Cells.Item(15, 8).Value = [double]"26.97695"
$ws.Cells.Item(15, 9).Value = [double]"0.08214718905254291"
$ws.Cells.Item(15, 10).Value = [double]"0.08214718905254291"
$ws.Cells.Item(15, 15).Value = [double]"0.4709815605157605"
$ws.Cells.Item(15, 16).Value = [double]"0.4709815605157605"
$ws.Cells.Item(15, 17).Value = [double]"520.633766258161"
$ws.Cells.Item(15, 18).Value = [double]"4685.70389632345"
$ws.Cells.Item(15, 19).Value = [double]"0.03868981129194985"
$ws.Cells.Item(15, 20).Value = [double]"0.03868981129194985"
$ws.Cells.Item(16, 7).Value = [double]"8.992316666666666"
$ws.Cells.Item(16, 8).Value = [double]"26.97695"
$ws.Cells.Item(16, 9).Value = [double]"0.08214718905254291"
$ws.Cells.Item(16, 10).Value = [double]"0.08214718905254291"
$ws.Cells.Item(16, 13).Value = [double]"61.10114166666667"
$ws.Cells.Item(16, 14).Value = [double]"183.303425"
$ws.Cells.Item(16, 15).Value = [double]"0.4970413158429724"
$ws.Cells.Item(16, 16).Value = [double]"0.4970413158429724"
$ws.Cells.Item(16, 17).Value = [double]"549.4408145615278"
$ws.Cells.Item(16, 18).Value = [double]"4944.96733105375"
$ws.Cells.Item(16, 19).Value = [double]"0.04083054693947735"
$ws.Cells.Item(16, 20).Value = [double]"0.04083054693947735"
$ws.Cells.Item(17, 7).Value = [double]"8.992316666666666"
$ws.Cells.Item(17, 8).Value = [double]"26.97695"
$ws.Cells.Item(17, 9).Value = [double]"0.08214718905254291"
$ws.Cells.Item(17, 10).Value = [double]"0.08214718905254291"
$ws.Cells.Item(17, 13).Value = [double]"0.1059033333333333"
$ws.Cells.Item(17, 14).Value = [double]"0.31771"
$ws.Cells.Item(17, 15).Value = [double]"0.0008614950672987739"
$ws.Cells.Item(17, 16).Value = [double]"0.0008614950672987739"
$ws.Cells.Item(17, 17).Value = [double]"0.9523163093888888"
$ws.Cells.Item(17, 18).Value = [double]"8.570846784499999"
$ws.Cells.Item(17, 19).Value = [double]"7.076939816122555E-05"
$ws.Cells.Item(17, 20).Value = [double]"7.076939816122555E-05"
